$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything
$ws.UsedRange.Delete()

# Re-type headers
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "User ID"
$ws.Range("C1").Value = "Année Fiscale"
$ws.Range("D1").Value = "Total Impôt"
$ws.Range("E1").Value = "Total Payé"
$ws.Range("F1").Value = "Statut"
$ws.Range("G1").Value = "Date Limite"
$ws.Range("H1").Value = "Moyen de Paiement"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2025
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 600
$ws.Range("F2").Value = "En cours"
$ws.Range("G2").Value = "'2025-02-05"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").Value = "Carte bancaire"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2025
$ws.Range("D3").Value = 100000
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = "En cours"
$ws.Range("G3").Value = "'2025-02-05"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = "cash"
